$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 195
$ws.Range("I18").Value = 195
$ws.Range("K18").Value = 195
$ws.Range("M18").Value = 89
$ws.Range("H28").Value = 320.25
$ws.Range("I28").Value = 144.4
$ws.Range("K28").Value = 144.4
$ws.Range("M28").Value = 340.6
$ws.Range("H62").Value = 3803.2632
$ws.Range("I62").Value = 2936.3635
$ws.Range("K62").Value = 2936.3635
$ws.Range("M62").Value = -2312.3635
$ws.Range("H64").Value = 3956.1875
$ws.Range("J64").Value = 4371.4287
$ws.Range("L64").Value = 4371.4287
$ws.Range("N64").Value = -4867.4287
$ws.Range("H65").Value = 3803.2632
$ws.Range("I65").Value = 2936.3635
$ws.Range("K65").Value = 14681.8175
$ws.Range("M65").Value = -11561.8175
$ws.Range("H67").Value = 3956.1875
$ws.Range("J67").Value = 4371.4287
$ws.Range("L67").Value = 4371.4287
$ws.Range("N67").Value = -6087.4287
$ws.Range("H74").Value = 4995.364
$ws.Range("I74").Value = 4499.5
$ws.Range("J74").Value = 5105.5557
$ws.Range("K74").Value = 4499.5
$ws.Range("L74").Value = 5105.5557
$ws.Range("M74").Value = -3563.5
$ws.Range("N74").Value = -6977.5557
$ws.Range("H76").Value = 3833.25
$ws.Range("I76").Value = 3666.5
$ws.Range("J76").Value = 4000
$ws.Range("K76").Value = 3666.5
$ws.Range("L76").Value = 4000
$ws.Range("M76").Value = -3351.5
$ws.Range("N76").Value = -4630
$ws.Range("H77").Value = 4995.364
$ws.Range("I77").Value = 4499.5
$ws.Range("J77").Value = 5105.5557
$ws.Range("K77").Value = 22497.5
$ws.Range("L77").Value = 25527.7785
$ws.Range("M77").Value = -17817.5
$ws.Range("N77").Value = -34887.7785
$ws.Range("H79").Value = 3833.25
$ws.Range("I79").Value = 3666.5
$ws.Range("J79").Value = 4000
$ws.Range("K79").Value = 3666.5
$ws.Range("L79").Value = 4000
$ws.Range("M79").Value = -2574.5
$ws.Range("N79").Value = -6184
$ws.Range("H86").Value = 7512.0557
$ws.Range("I86").Value = 2204.0908
$ws.Range("J86").Value = 15853.143
$ws.Range("K86").Value = 2204.0908
$ws.Range("L86").Value = 15853.143
$ws.Range("M86").Value = -1081.0908
$ws.Range("N86").Value = -18099.143
$ws.Range("H89").Value = 7512.0557
$ws.Range("I89").Value = 2204.0908
$ws.Range("J89").Value = 15853.143
$ws.Range("K89").Value = 11020.454
$ws.Range("L89").Value = 79265.715
$ws.Range("M89").Value = -5404.454
$ws.Range("N89").Value = -90497.715
$ws.Range("H92").Value = 100000780
$ws.Range("I92").Value = 111111970
$ws.Range("K92").Value = 111111970
$ws.Range("M92").Value = -111110722
$ws.Range("H112").Value = 1140
$ws.Range("J112").Value = 1140
$ws.Range("L112").Value = 3420
$ws.Range("N112").Value = -5636
$ws.Range("H113").Value = 31253708
$ws.Range("I113").Value = 71431780
$ws.Range("J113").Value = 4096.9443
$ws.Range("K113").Value = 71431780
$ws.Range("L113").Value = 4096.9443
$ws.Range("M113").Value = -71428526
$ws.Range("N113").Value = -10604.9443
$ws.Range("H129").Value = 1124.4546
$ws.Range("I129").Value = 465.7143
$ws.Range("J129").Value = 1220.5209
$ws.Range("K129").Value = 1397.1429
$ws.Range("L129").Value = 3661.5627
$ws.Range("M129").Value = 3602.8571
$ws.Range("N129").Value = -13661.5627
$ws.Range("H138").Value = 2374.8147
$ws.Range("I138").Value = 3133.3333
$ws.Range("J138").Value = 2280
$ws.Range("K138").Value = 9399.999899999999
$ws.Range("L138").Value = 6840
$ws.Range("M138").Value = -4259.999899999999
$ws.Range("N138").Value = -17120

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5052.07
$ws.Range("I32").Value = 5258.5
$ws.Range("K32").Value = 5258.5
$ws.Range("M32").Value = -4971.5
$ws.Range("H97").Value = 1469.5555
$ws.Range("I97").Value = 1563.7142
$ws.Range("J97").Value = 1140
$ws.Range("K97").Value = 1563.7142
$ws.Range("L97").Value = 1140
$ws.Range("M97").Value = -1067.7142
$ws.Range("N97").Value = -2132
$ws.Range("H110").Value = 702.625
$ws.Range("I110").Value = 688.7143
$ws.Range("K110").Value = 688.7143
$ws.Range("M110").Value = 1356.2857
$ws.Range("H122").Value = 2952.1177
$ws.Range("I122").Value = 2507.3333
$ws.Range("K122").Value = 7521.999899999999
$ws.Range("M122").Value = -5071.999899999999
$ws.Range("H132").Value = 15157.622
$ws.Range("I132").Value = 1543.2258
$ws.Range("J132").Value = 85498.664
$ws.Range("K132").Value = 4629.6774
$ws.Range("L132").Value = 256495.992
$ws.Range("M132").Value = -2099.6774
$ws.Range("N132").Value = -261555.992

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1367.8125
$ws.Range("I107").Value = 1355.4445
$ws.Range("J107").Value = 1383.7142
$ws.Range("K107").Value = 1355.4445
$ws.Range("L107").Value = 1383.7142
$ws.Range("M107").Value = 564.5554999999999
$ws.Range("N107").Value = -5223.7142

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1166
$ws.Range("I16").Value = 1177.7778
$ws.Range("J16").Value = 1139.5
$ws.Range("K16").Value = 1177.7778
$ws.Range("L16").Value = 1139.5
$ws.Range("M16").Value = -890.7778000000001
$ws.Range("N16").Value = -1713.5
$ws.Range("H22").Value = 240.55556
$ws.Range("I22").Value = 130.5
$ws.Range("J22").Value = 328.6
$ws.Range("K22").Value = 130.5
$ws.Range("L22").Value = 328.6
$ws.Range("M22").Value = 219.5
$ws.Range("N22").Value = -1028.6
$ws.Range("H58").Value = 13392.05
$ws.Range("I58").Value = 817.9677
$ws.Range("J58").Value = 56702.777
$ws.Range("K58").Value = 817.9677
$ws.Range("L58").Value = 56702.777
$ws.Range("M58").Value = -614.9677
$ws.Range("N58").Value = -57108.777
$ws.Range("H113").Value = 1166
$ws.Range("I113").Value = 1177.7778
$ws.Range("J113").Value = 1139.5
$ws.Range("K113").Value = 1177.7778
$ws.Range("L113").Value = 1139.5
$ws.Range("M113").Value = 992.2221999999999
$ws.Range("N113").Value = -5479.5
$ws.Range("H136").Value = 13392.05
$ws.Range("I136").Value = 817.9677
$ws.Range("J136").Value = 56702.777
$ws.Range("K136").Value = 2453.9031
$ws.Range("L136").Value = 170108.331
$ws.Range("M136").Value = 96.09690000000001
$ws.Range("N136").Value = -175208.331

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2797
$ws.Range("I5").Value = 5063.7144
$ws.Range("J5").Value = 813.625
$ws.Range("K5").Value = 15191.1432
$ws.Range("L5").Value = 2440.875
$ws.Range("M5").Value = -15079.1432
$ws.Range("N5").Value = -2664.875
$ws.Range("H122").Value = 526.7222
$ws.Range("I122").Value = 231.91667
$ws.Range("J122").Value = 1116.3334
$ws.Range("K122").Value = 2087.25003
$ws.Range("L122").Value = 10047.0006
$ws.Range("M122").Value = 362.7499699999998
$ws.Range("N122").Value = -14947.0006
$ws.Range("H131").Value = 795.51
$ws.Range("J131").Value = 803.31915
$ws.Range("L131").Value = 2409.95745
$ws.Range("N131").Value = -12489.95745
$ws.Range("H135").Value = 2797
$ws.Range("I135").Value = 5063.7144
$ws.Range("J135").Value = 813.625
$ws.Range("K135").Value = 45573.4296
$ws.Range("L135").Value = 7322.625
$ws.Range("M135").Value = -43038.4296
$ws.Range("N135").Value = -12392.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4247.846
$ws.Range("I80").Value = 3551.6667
$ws.Range("J80").Value = 4844.5713
$ws.Range("K80").Value = 3551.6667
$ws.Range("L80").Value = 4844.5713
$ws.Range("M80").Value = -2553.6667
$ws.Range("N80").Value = -6840.5713
$ws.Range("H83").Value = 4247.846
$ws.Range("I83").Value = 3551.6667
$ws.Range("J83").Value = 4844.5713
$ws.Range("K83").Value = 17758.3335
$ws.Range("L83").Value = 24222.8565
$ws.Range("M83").Value = -12766.3335
$ws.Range("N83").Value = -34206.85649999999
$ws.Range("H95").Value = 4211
$ws.Range("J95").Value = 4211
$ws.Range("L95").Value = 4211
$ws.Range("N95").Value = -9703
$ws.Range("H132").Value = 47816.582
$ws.Range("I132").Value = 7089
$ws.Range("K132").Value = 21267
$ws.Range("M132").Value = -18737

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2618.4
$ws.Range("I82").Value = 3580
$ws.Range("J82").Value = 2137.6
$ws.Range("K82").Value = 3580
$ws.Range("L82").Value = 2137.6
$ws.Range("M82").Value = -3219
$ws.Range("N82").Value = -2859.6
$ws.Range("H85").Value = 2618.4
$ws.Range("I85").Value = 3580
$ws.Range("J85").Value = 2137.6
$ws.Range("K85").Value = 3580
$ws.Range("L85").Value = 2137.6
$ws.Range("M85").Value = -2332
$ws.Range("N85").Value = -4633.6
$ws.Range("H132").Value = 1918.2727
$ws.Range("I132").Value = 1267.0667
$ws.Range("J132").Value = 3313.7144
$ws.Range("K132").Value = 3801.2001
$ws.Range("L132").Value = 9941.143199999999
$ws.Range("M132").Value = -1271.2001
$ws.Range("N132").Value = -15001.1432
$ws.Range("H136").Value = 1346.8667
$ws.Range("I136").Value = 1260.4
$ws.Range("J136").Value = 1519.8
$ws.Range("K136").Value = 3781.2
$ws.Range("L136").Value = 4559.4
$ws.Range("M136").Value = -1231.2
$ws.Range("N136").Value = -9659.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 71429670
$ws.Range("I81").Value = 1454.5
$ws.Range("J81").Value = 250000210
$ws.Range("K81").Value = 2909
$ws.Range("L81").Value = 500000420
$ws.Range("M81").Value = -1848
$ws.Range("N81").Value = -500002542
$ws.Range("H84").Value = 71429670
$ws.Range("I84").Value = 1454.5
$ws.Range("J84").Value = 250000210
$ws.Range("K84").Value = 14545
$ws.Range("L84").Value = 2500002100
$ws.Range("M84").Value = -9241
$ws.Range("N84").Value = -2500012708
$ws.Range("H132").Value = 1019.5897
$ws.Range("I132").Value = 723.5925999999999
$ws.Range("J132").Value = 1685.5834
$ws.Range("K132").Value = 2170.7778
$ws.Range("L132").Value = 5056.7502
$ws.Range("M132").Value = 359.2222000000002
$ws.Range("N132").Value = -10116.7502
$ws.Range("H136").Value = 35715916
$ws.Range("I136").Value = 43479750
$ws.Range("K136").Value = 130439250
$ws.Range("M136").Value = -130436700
